$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the taxon-related data between row 26 and row 28
# (columns A, B, E, F, G, H, Q, R), while leaving the rest of the
# rows (C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AI, AT,
# AW, AX, AY) untouched since they already match between the two rows.

$cols = @("A", "B", "E", "F", "G", "H", "Q", "R")

foreach ($col in $cols) {
    $addr26 = "$col" + "26"
    $addr28 = "$col" + "28"

    $val26 = $ws.Range($addr26).Value2
    $val28 = $ws.Range($addr28).Value2

    $ws.Range($addr26).Value2 = $val28
    $ws.Range($addr28).Value2 = $val26
}
